# Refresh the cryptos list with the latest scraped price/volume figures.
# Note: some Price (column D) values look like plain numbers (e.g. "484.32"),
# but the source data stores them as text. A leading apostrophe forces Excel
# to keep them as text instead of coercing to a Double (which would drop
# trailing zeros, e.g. "10.60" -> 10.6).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').Value = '68.479.86'
$ws.Range('E2').Value = '  +1.42%  '
$ws.Range('D3').Value = '3.931.58'
$ws.Range('E3').Value = '  +1.91%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '''484.32'
$ws.Range('E5').Value = '  +5.53%  '
$ws.Range('D6').Value = '''148.56'
$ws.Range('E6').Value = '  +1.77%  '
$ws.Range('E7').Value = '  -0.66%  '
$ws.Range('D8').Value = '''0.998'
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('D9').Value = '''0.724'
$ws.Range('E9').Value = '  -3.23%  '
$ws.Range('E10').Value = '  +8.26%  '
$ws.Range('E11').Value = '  +11.50%  '
$ws.Range('D12').Value = '''42.53'
$ws.Range('E12').Value = '  -2.87%  '
$ws.Range('D13').Value = '''10.60'
$ws.Range('E13').Value = '  +1.91%  '
$ws.Range('D14').Value = '4.558.36'
$ws.Range('E14').Value = '  +2.21%  '
$ws.Range('D15').Value = '3.985.46'
$ws.Range('E15').Value = '  +4.34%  '
$ws.Range('D16').Value = '''14.63'
$ws.Range('E16').Value = '  -0.62%  '
$ws.Range('D18').Value = '''19.76'
$ws.Range('E18').Value = '  -1.32%  '
$ws.Range('E19').Value = '  -2.40%  '
$ws.Range('D20').Value = '68.571.71'
$ws.Range('E20').Value = '  +1.38%  '
$ws.Range('D21').Value = '''432.08'
$ws.Range('E21').Value = '  +1.25%  '
$ws.Range('D22').Value = '''3.35'
$ws.Range('E22').Value = '  +2.28%  '
$ws.Range('D23').Value = '''14.50'
$ws.Range('E23').Value = '  -2.42%  '
$ws.Range('D24').Value = '''87.06'
$ws.Range('D25').Value = '''11.31'
$ws.Range('E25').Value = '  +13.30%  '
$ws.Range('E26').Value = '  +0.77%  '
$ws.Range('B27').Value = 'RenderToken'
$ws.Range('C27').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D27').Value = '''10.49'
$ws.Range('E27').Value = '  -0.50%  '
$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').Value = '''38.18'
$ws.Range('E28').Value = '  +1.68%  '
$ws.Range('E29').Value = '  +6.78%  '
$ws.Range('D30').Value = '''718.52'
$ws.Range('E30').Value = '  -4.59%  '
$ws.Range('E31').Value = '  -3.93%  '
$ws.Range('E32').Value = '  -4.65%  '
$ws.Range('E33').Value = '  +3.34%  '
$ws.Range('D34').Value = '0.0₃0888'
$ws.Range('E34').Value = '  +31.12%  '
$ws.Range('D35').Value = '''41.70'
$ws.Range('E35').Value = '  -5.09%  '
$ws.Range('D36').Value = '''59.08'
$ws.Range('E36').Value = '  +2.71%  '
$ws.Range('E37').Value = '  -7.12%  '
$ws.Range('D38').Value = '''5.51'
$ws.Range('E38').Value = '  -0.91%  '
$ws.Range('E39').Value = '  -0.15%  '
$ws.Range('D40').Value = '''2.86'
$ws.Range('E40').Value = '  +7.80%  '
$ws.Range('D41').Value = '''0.0468'
$ws.Range('E41').Value = '  -1.95%  '
$ws.Range('E42').Value = '  +10.05%  '
$ws.Range('E43').Value = '  +1.57%  '
$ws.Range('E44').Value = '  -3.79%  '
$ws.Range('D45').Value = '''0.141'
$ws.Range('E45').Value = '  -0.01%  '
$ws.Range('E46').Value = '  +0.02%  '
$ws.Range('B47').Value = 'ARBITRUM'
$ws.Range('C47').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D47').Value = '''2.17'
$ws.Range('E47').Value = '  +1.87%  '
$ws.Range('B48').Value = 'LidoDAOToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D48').Value = '''3.41'
$ws.Range('E48').Value = '  -0.99%  '
$ws.Range('E49').Value = '  -0.09%  '
$ws.Range('D50').Value = '''146.80'
$ws.Range('E50').Value = '  +1.60%  '
$ws.Range('D51').Value = '''2.84'
$ws.Range('E51').Value = '  -1.63%  '
